$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs in the "Focus on updating the documentation..."
#    bullet back into a single run, removing the _GoBack bookmark that used
#    to sit between them.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$d.Content.Find.Execute(
    "Focus on updating the documentation daily rather than sporadically and upload each new version to GitHub rather than working on it locally as that does not provide evidence of our daily work.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Focus on updating the documentation daily rather than sporadically and upload each new version to GitHub rather than working on it locally as that does not provide evidence of our daily work.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) Close out the "Burndown chart" bullet with a trailing period, then add a
#    brand new bullet about running better scrum meetings (Rauf's
#    suggestion). Move the _GoBack bookmark to the end of this new bullet.
# ---------------------------------------------------------------------------
$burndownPara = $d.Paragraphs(9)
$burndownPara.Range.InsertAfter(".")

$burndownPara = $d.Paragraphs(9)
$burndownPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs(10)
# Append a placeholder character so the bookmark insertion point below isn't
# the very last character position in the paragraph (that boundary case
# misplaces the bookmark), then strip the placeholder back out afterwards.
$newPara.Range.InsertAfter("Conduct our scrum meetings better by making it less of a report and more of a discussion, more concise, update tasks while we are talking.X")

$newPara = $d.Paragraphs(10)
$bookmarkPos = $newPara.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholder = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$placeholder.Delete()
